$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new changelog entry in column A, row 12
$ws.Range("A12").Value = "Added support for Polytone"

# Update selection to reflect the next empty row (A13), matching Excel's
# behavior of moving the active cell after the last used row.
$ws.Range("A13").Select()
